$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.385.17"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.60%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.190.82"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.96%  "

$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "250.34"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.84%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.610"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.44%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "73.46"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.53%  "

$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.584"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.40%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.93"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.36%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0913"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.74%  "

$ws.Range("E12").Value = "  +1.12%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.76"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.84%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.521.49"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.87%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.30"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.02%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.180.50"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.43%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.777"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.40%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.311.06"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.76%  "

$ws.Range("E19").Value = "  +0.87%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.91"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.63%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.87"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.84%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "227.55"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.37%  "

$ws.Range("E23").Value = "  +9.24%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.37"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -7.25%  "

$ws.Range("E25").Value = "  -0.07%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.62"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.84%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.37"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.56%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "39.44"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +21.67%  "

$ws.Range("E29").Value = "  +2.54%  "

$ws.Range("E30").Value = "  +0.55%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "169.46"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.31%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.01"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.99%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0792"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.70%  "

$ws.Range("E34").Value = "  +0.78%  "

$ws.Range("E35").Value = "  -0.45%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.106"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.52%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.32"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.43%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0324"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +7.40%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "12.00"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.34%  "

$ws.Range("E40").Value = "  +0.15%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.198"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.35%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.23"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.21%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "58.47"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.26%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.73"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.09%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.46"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.68%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.472"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +19.17%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0971"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.27%  "

$ws.Range("E48").Value = "  +9.59%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.09"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.19%  "

$ws.Range("E50").Value = "  +0.32%  "

$ws.Range("E51").Value = "  +1.81%  "
